$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2,4)
$c.NumberFormat = "@"
$c.Value = '30.625.99'
$c.Style = "Normal"
$c = $ws.Cells.Item(3,4)
$c.NumberFormat = "@"
$c.Value = '2.112.41'
$c.Style = "Normal"
$c = $ws.Cells.Item(3,5)
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(4,4)
$c.NumberFormat = "@"
$c.Value = '1.014'
$c.Style = "Normal"
$c = $ws.Cells.Item(4,5)
$c.NumberFormat = "@"
$c.Value = '  +1.18%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value = '338.66'
$c.Style = "Normal"
$c = $ws.Cells.Item(5,5)
$c.NumberFormat = "@"
$c.Value = '  +1.41%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value = '1.012'
$c.Style = "Normal"
$c = $ws.Cells.Item(6,5)
$c.NumberFormat = "@"
$c.Value = '  +1.07%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(7,4)
$c.NumberFormat = "@"
$c.Value = '0.5250'
$c.Style = "Normal"
$c = $ws.Cells.Item(7,5)
$c.NumberFormat = "@"
$c.Value = '  -0.31%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(8,5)
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(9,4)
$c.NumberFormat = "@"
$c.Value = '53.46'
$c.Style = "Normal"
$c = $ws.Cells.Item(9,5)
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(10,4)
$c.NumberFormat = "@"
$c.Value = '0.09027'
$c.Style = "Normal"
$c = $ws.Cells.Item(10,5)
$c.NumberFormat = "@"
$c.Value = '  +0.19%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(11,4)
$c.NumberFormat = "@"
$c.Value = '1.168'
$c.Style = "Normal"
$c = $ws.Cells.Item(11,5)
$c.NumberFormat = "@"
$c.Value = '  -1.03%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(12,4)
$c.NumberFormat = "@"
$c.Value = '24.35'
$c.Style = "Normal"
$c = $ws.Cells.Item(12,5)
$c.NumberFormat = "@"
$c.Value = '  -0.44%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(13,4)
$c.NumberFormat = "@"
$c.Value = '2.126.11'
$c.Style = "Normal"
$c = $ws.Cells.Item(13,5)
$c.NumberFormat = "@"
$c.Value = '  +1.26%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(14,4)
$c.NumberFormat = "@"
$c.Value = '6.777'
$c.Style = "Normal"
$c = $ws.Cells.Item(14,5)
$c.NumberFormat = "@"
$c.Value = '  -0.22%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(15,4)
$c.NumberFormat = "@"
$c.Value = '8.057'
$c.Style = "Normal"
$c = $ws.Cells.Item(15,5)
$c.NumberFormat = "@"
$c.Value = '  +3.13%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(16,4)
$c.NumberFormat = "@"
$c.Value = '97.86'
$c.Style = "Normal"
$c = $ws.Cells.Item(17,4)
$c.NumberFormat = "@"
$c.Value = '0.00001162'
$c.Style = "Normal"
$c = $ws.Cells.Item(17,5)
$c.NumberFormat = "@"
$c.Value = '  +2.61%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(18,4)
$c.NumberFormat = "@"
$c.Value = '1.014'
$c.Style = "Normal"
$c = $ws.Cells.Item(18,5)
$c.NumberFormat = "@"
$c.Value = '  +1.05%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(19,4)
$c.NumberFormat = "@"
$c.Value = '0.06704'
$c.Style = "Normal"
$c = $ws.Cells.Item(19,5)
$c.NumberFormat = "@"
$c.Value = '  +1.20%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(20,4)
$c.NumberFormat = "@"
$c.Value = '19.31'
$c.Style = "Normal"
$c = $ws.Cells.Item(20,5)
$c.NumberFormat = "@"
$c.Value = '  -1.05%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(21,4)
$c.NumberFormat = "@"
$c.Value = '1.012'
$c.Style = "Normal"
$c = $ws.Cells.Item(21,5)
$c.NumberFormat = "@"
$c.Value = '  +1.04%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(22,4)
$c.NumberFormat = "@"
$c.Value = '6.321'
$c.Style = "Normal"
$c = $ws.Cells.Item(22,5)
$c.NumberFormat = "@"
$c.Value = '  -0.10%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(23,4)
$c.NumberFormat = "@"
$c.Value = '30.728.63'
$c.Style = "Normal"
$c = $ws.Cells.Item(23,5)
$c.NumberFormat = "@"
$c.Value = '  +0.64%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(24,4)
$c.NumberFormat = "@"
$c.Value = '12.79'
$c.Style = "Normal"
$c = $ws.Cells.Item(24,5)
$c.NumberFormat = "@"
$c.Value = '  +2.99%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(25,4)
$c.NumberFormat = "@"
$c.Value = '2.377'
$c.Style = "Normal"
$c = $ws.Cells.Item(25,5)
$c.NumberFormat = "@"
$c.Value = '  +1.02%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(26,4)
$c.NumberFormat = "@"
$c.Value = '2.373.85'
$c.Style = "Normal"
$c = $ws.Cells.Item(26,5)
$c.NumberFormat = "@"
$c.Value = '  +1.19%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(27,4)
$c.NumberFormat = "@"
$c.Value = '22.31'
$c.Style = "Normal"
$c = $ws.Cells.Item(27,5)
$c.NumberFormat = "@"
$c.Value = '  -0.29%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(28,5)
$c.NumberFormat = "@"
$c.Value = '  +1.04%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(29,4)
$c.NumberFormat = "@"
$c.Value = '2.535'
$c.Style = "Normal"
$c = $ws.Cells.Item(29,5)
$c.NumberFormat = "@"
$c.Value = '  -1.90%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(30,4)
$c.NumberFormat = "@"
$c.Value = '134.71'
$c.Style = "Normal"
$c = $ws.Cells.Item(30,5)
$c.NumberFormat = "@"
$c.Value = '  +1.29%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(31,4)
$c.NumberFormat = "@"
$c.Value = '1.192'
$c.Style = "Normal"
$c = $ws.Cells.Item(31,5)
$c.NumberFormat = "@"
$c.Value = '  -0.67%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(32,5)
$c.NumberFormat = "@"
$c.Value = '  -0.25%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(33,4)
$c.NumberFormat = "@"
$c.Value = '6.365'
$c.Style = "Normal"
$c = $ws.Cells.Item(33,5)
$c.NumberFormat = "@"
$c.Value = '  +3.24%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(34,5)
$c.NumberFormat = "@"
$c.Value = '  -2.10%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(35,4)
$c.NumberFormat = "@"
$c.Value = '3.947'
$c.Style = "Normal"
$c = $ws.Cells.Item(35,5)
$c.NumberFormat = "@"
$c.Value = '  +0.55%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(36,4)
$c.NumberFormat = "@"
$c.Value = '10.29'
$c.Style = "Normal"
$c = $ws.Cells.Item(36,5)
$c.NumberFormat = "@"
$c.Value = '  -2.66%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(37,4)
$c.NumberFormat = "@"
$c.Value = '5.882'
$c.Style = "Normal"
$c = $ws.Cells.Item(37,5)
$c.NumberFormat = "@"
$c.Value = '  +5.22%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(38,4)
$c.NumberFormat = "@"
$c.Value = '0.02649'
$c.Style = "Normal"
$c = $ws.Cells.Item(38,5)
$c.NumberFormat = "@"
$c.Value = '  +2.57%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(39,4)
$c.NumberFormat = "@"
$c.Value = '0.06827'
$c.Style = "Normal"
$c = $ws.Cells.Item(39,5)
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value = '0.2315'
$c.Style = "Normal"
$c = $ws.Cells.Item(40,5)
$c.NumberFormat = "@"
$c.Value = '  +0.37%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(41,4)
$c.NumberFormat = "@"
$c.Value = '12.58'
$c.Style = "Normal"
$c = $ws.Cells.Item(41,5)
$c.NumberFormat = "@"
$c.Value = '  -1.57%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(42,4)
$c.NumberFormat = "@"
$c.Value = '0.6869'
$c.Style = "Normal"
$c = $ws.Cells.Item(42,5)
$c.NumberFormat = "@"
$c.Value = '  -0.83%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(43,4)
$c.NumberFormat = "@"
$c.Value = '1.258'
$c.Style = "Normal"
$c = $ws.Cells.Item(43,5)
$c.NumberFormat = "@"
$c.Value = '  +0.64%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(44,4)
$c.NumberFormat = "@"
$c.Value = '14.93'
$c.Style = "Normal"
$c = $ws.Cells.Item(44,5)
$c.NumberFormat = "@"
$c.Value = '  +5.82%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(45,4)
$c.NumberFormat = "@"
$c.Value = '0.6421'
$c.Style = "Normal"
$c = $ws.Cells.Item(45,5)
$c.NumberFormat = "@"
$c.Value = '  +0.42%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(46,4)
$c.NumberFormat = "@"
$c.Value = '2.310'
$c.Style = "Normal"
$c = $ws.Cells.Item(46,5)
$c.NumberFormat = "@"
$c.Value = '  -2.09%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(47,5)
$c.NumberFormat = "@"
$c.Value = '  +11.23%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(48,4)
$c.NumberFormat = "@"
$c.Value = '3.701'
$c.Style = "Normal"
$c = $ws.Cells.Item(48,5)
$c.NumberFormat = "@"
$c.Value = '  +1.12%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(49,5)
$c.NumberFormat = "@"
$c.Value = '  +0.34%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(50,4)
$c.NumberFormat = "@"
$c.Value = '82.78'
$c.Style = "Normal"
$c = $ws.Cells.Item(50,5)
$c.NumberFormat = "@"
$c.Value = '  -1.12%  '
$c.Style = "Normal"
$c = $ws.Cells.Item(51,4)
$c.NumberFormat = "@"
$c.Value = '0.07306'
$c.Style = "Normal"
